# Restore revision 34f6a51b... (SAVE) — single-cell value fix on the
# "Rules" sheet: rule R30's "From" threshold (C10) is corrected from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
